# Applies: "Added start parse at row index. Replaced any value type in
# ignore field on only bool type (available value: empty, "true", "false",
# "1", "0"). If type is incorrect printed warning"
#
# Concretely, in the example workbook this means the "ignore" column (B)
# on sheet "Lvl0" now holds boolean values instead of arbitrary text/number,
# two new example rows (9/10) were added showing the new "ignoredFeature"
# strings, and the user ended up with sheet "Lvl0" active (cell F8
# selected) instead of sheet "Lvl2".

$wb = $excel.ActiveWorkbook

$lvl0 = $wb.Worksheets.Item("Lvl0")

# --- New example rows 9 & 10: ignore column now boolean, with a note in C ---
# Row 10 previously held the string "qw" in the ignore column (B) -- that
# was an invalid/arbitrary value for a field that is now strictly boolean.
# (Row 10 is filled in before row 9 so the shared-string table ends up in
# the same order as the authored workbook.)
$lvl0.Cells.Item(10, 2).Value = $true
$lvl0.Cells.Item(10, 3).Value = "ignoredFeature2"

$lvl0.Cells.Item(9, 2).Value = $true
$lvl0.Cells.Item(9, 3).Value = "ignoredFeature1"

# --- Row 25/26: ignore column value 123 (number) is no longer a valid
# bool-ish value, replaced with proper boolean false; row 26 also gets an
# explicit boolean false in the ignore column. ---
$lvl0.Cells.Item(25, 2).Value = $false
$lvl0.Cells.Item(26, 2).Value = $false

# Column C got wider to fit the new "ignoredFeature1"/"ignoredFeature2" text.
$lvl0.Columns.Item(3).ColumnWidth = 15.02

# --- View state: user ends up on the "Lvl0" tab with F8 selected (rather
# than "Lvl2" with tabSelected) ---
$lvl0.Activate() | Out-Null
$lvl0.Range("F8").Select() | Out-Null
